$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 data
$ws.Range("A6").Value2 = "1557. Minimum Number of Vertices to Reach All Nodes"

$ws.Range("B6").Value2 = "Medium"
$ws.Range("B6").Interior.Color = $ws.Range("B2").Interior.Color

$ws.Range("C6").Value2 = "Graph Theory"

$ws.Range("D6").Value2 = "Basically find all the nodes with in-degree = 0. Because you can't get to these nodes from anywhere, and once you have all of these you can visit the next level of nodes, and so on. We don't need to do graph traversal, but just look at the vertices with a seen[] map."

$ws.Hyperlinks.Add($ws.Range("E6"), "https://leetcode.com/problems/minimum-number-of-vertices-to-reach-all-nodes/solutions/805685/java-c-python-nodes-with-no-in-degree/?envType=study-plan-v2&envId=graph-theory ", "", "", "https://leetcode.com/problems/minimum-number-of-vertices-to-reach-all-nodes/solutions/805685/java-c-python-nodes-with-no-in-degree/?envType=study-plan-v2&envId=graph-theory ") | Out-Null
$ws.Range("E6").Style = $ws.Range("E2").Style

# Update selection to match the saved workbook state
$ws.Range("D7").Select()
